# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on Sheet1, matching the latest GitHub Actions data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D (Price) to be stored as text so values like "1.000" or
# "307.80" are not re-interpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row = 2;  D = "27.161.59";     E = "  +0.72%  " }
    @{ Row = 3;  D = "1.907.46";      E = "  +1.75%  " }
    @{ Row = 4;  D = "1.000";         E = "  +0.03%  " }
    @{ Row = 5;  D = "307.80";        E = "  +0.77%  " }
    @{ Row = 6;  D = "1.000";         E = "  -0.01%  " }
    @{ Row = 7;  D = "0.5242";        E = "  +3.03%  " }
    @{ Row = 8;  D = "0.3780";        E = "  +2.99%  " }
    @{ Row = 9;  D = "0.07269";       E = "  +0.81%  " }
    @{ Row = 10; D = "21.29";         E = "  +2.46%  " }
    @{ Row = 11; D = "0.8980";        E = "  +0.27%  " }
    @{ Row = 12; D = "0.07678";       E = "  +2.07%  " }
    @{ Row = 13; D = "1.884.15";      E = "  +0.42%  " }
    @{ Row = 14; D = "94.77";         E = "  -0.42%  " }
    @{ Row = 15; D = "5.249";         E = "  +0.02%  " }
    @{ Row = 16; D = "1.001";         E = "  +0.07%  " }
    @{ Row = 17; D = "0.000008548";   E = "  +0.13%  " }
    @{ Row = 18; D = "14.56";         E = "  +2.10%  " }
    @{ Row = 19; D = "1.0000";        E = "  +0.01%  " }
    @{ Row = 20; D = "27.211.87";     E = "  +0.76%  " }
    @{ Row = 21; D = "5.088";         E = "  +1.23%  " }
    @{ Row = 22; D = "2.136.95";      E = "  +0.80%  " }
    @{ Row = 23; D = "10.64";         E = "  +2.37%  " }
    @{ Row = 24; D = "6.448";         E = "  +0.65%  " }
    @{ Row = 25; D = "2.316";         E = "  +10.58%  " }
    @{ Row = 26; D = "145.80";        E = "  -1.87%  " }
    @{ Row = 27; D = "18.16";         E = "  +1.37%  " }
    @{ Row = 28; D = "1.731";         E = "  -3.39%  " }
    @{ Row = 29; D = "114.87";        E = "  +1.19%  " }
    @{ Row = 30; D = "4.966";         E = "  +4.81%  " }
    @{ Row = 31; D = "4.814";         E = "  +1.81%  " }
    @{ Row = 32; D = "0.09211";       E = "  +0.56%  " }
    @{ Row = 33; D = "0.05070";       E = "  -0.67%  " }
    @{ Row = 34; D = $null;           E = "  +7.29%  " }
    @{ Row = 35; D = "0.7803";        E = "  +3.95%  " }
    @{ Row = 36; D = $null;           E = "  +0.53%  " }
    @{ Row = 37; D = "3.312";         E = "  +2.24%  " }
    @{ Row = 38; D = $null;           E = "  +2.81%  " }
    @{ Row = 39; D = "0.5686";        E = "  +0.76%  " }
    @{ Row = 40; D = $null;           E = "  -0.51%  " }
    @{ Row = 41; D = "1.075";         E = "  -0.32%  " }
    @{ Row = 42; D = "9.040";         E = "  +5.26%  " }
    @{ Row = 43; D = "6.638";         E = "  -0.12%  " }
    @{ Row = 44; D = "118.72";        E = "  +2.92%  " }
    @{ Row = 45; D = "0.1524";        E = "  +2.99%  " }
    @{ Row = 46; D = "0.4858";        E = "  +2.17%  " }
    @{ Row = 47; D = "10.22";         E = "  +0.56%  " }
    @{ Row = 48; D = "1.000";         E = "  +0.04%  " }
    @{ Row = 49; D = "1.605";         E = "  +2.02%  " }
    @{ Row = 50; D = "37.48";         E = "  +1.31%  " }
    @{ Row = 51; D = "64.27";         E = "  +1.59%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
